$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ideal" numbers
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Update the active selection to match the saved view state
$ws.Range("E4").Select()

# Update the workbook window position/size
$excel.Left = 14520
$excel.Top = 0
$excel.Width = 14280
$excel.Height = 18000
